# Apply updated coin symbol list values scraped on 2022-12-13.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'259.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.151"
$ws.Range("D4").Style = "Normal"
$ws.Range("D7").Value = "'6.497"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.331"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8222"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.01327"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.1602"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.08090"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "'0.03188"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.09219"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.775"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.001644"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'0.006482"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.006131"
$ws.Range("D20").Style = "Normal"
$ws.Range("D24").Value = "'2.269"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3316"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1244"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = "'0.0002719"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.04585"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007019"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1115"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = "'0.003378"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = "'0.01121"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006053"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.0009922"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.8043"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Value = "'0.01243"
$ws.Range("D51").Style = "Normal"
